$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.619.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.982.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.72%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.635"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.83%  "
$ws.Range("E7").Value = "  +7.67%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("E12").Value = "  -2.80%  "
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.271.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.981.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "35.557.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0847"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.28%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +21.03%  "
$ws.Range("E26").Value = "  -4.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.50%  "
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.11%  "
$ws.Range("E32").Value = "  -6.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0954"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +15.16%  "
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.65%  "
$ws.Range("E36").Value = "  -3.55%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.16%  "
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.70%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0899"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.377.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "47.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("E51").Value = "  -0.10%  "

Write-Host "Applied crypto list update"
